# Resize/reposition the KPI-number textboxes on slides 4 and 5 so the
# (now larger) numbers keep fitting their boxes, and bump the font size
# / spacing on two of the slide-5 boxes to match.
#
# NOTE: Shape.Left/Top/Width/Height are expressed in points in this COM
# model, but are stored internally in EMU (1 pt = 12700 EMU) by flooring
# the point value (no standard rounding). To land exactly on a target
# EMU value we use the midpoint of the EMU's unit point-interval:
#   pt = (targetEmu + 0.5) / 12700

$p = $ppt.ActivePresentation

function Set-ShapeBox($Slide, $Name, $LeftEmu, $TopEmu, $WidthEmu, $HeightEmu) {
    $sh = $Slide.Shapes.Item($Name)
    $sh.Left = ($LeftEmu + 0.5) / 12700
    $sh.Top = ($TopEmu + 0.5) / 12700
    $sh.Width = ($WidthEmu + 0.5) / 12700
    $sh.Height = ($HeightEmu + 0.5) / 12700
}

# ---- Slide 4 ----
$s4 = $p.Slides.Item(4)

Set-ShapeBox $s4 "TextBox 94" 4983639 3417534 667986 491855
Set-ShapeBox $s4 "TextBox 95" 8042400 3431562 694808 491855
Set-ShapeBox $s4 "TextBox 96" 10938243 3417534 759055 491855
Set-ShapeBox $s4 "TextBox 97" 13897573 3417534 788751 491855

# ---- Slide 5 ----
$s5 = $p.Slides.Item(5)

Set-ShapeBox $s5 "TextBox 111" 3229800 3340381 1177460 480194
Set-ShapeBox $s5 "TextBox 112" 7166866 3340381 826107 480194
Set-ShapeBox $s5 "TextBox 113" 11106005 3325021 1008187 491855
Set-ShapeBox $s5 "TextBox 114" 15208594 3325021 1051337 491855

# TextBox 111 and TextBox 112 also shrink their line spacing / font size
# / character spacing slightly so the bigger box still looks right.
foreach ($name in @("TextBox 111", "TextBox 112")) {
    $sh = $s5.Shapes.Item($name)
    $tr = $sh.TextFrame.TextRange
    $tr.ParagraphFormat.SpaceWithin = 39.82
    $tr.Font.Size = 28.44
    $tr.Font.Spacing = -0.31
}
